$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.194828579899138
$ws.Range("C2").Value = 0.2170881049805047
$ws.Range("D2").Value = 0.02325209133014994
$ws.Range("E2").Value = 0.09891205464730568
$ws.Range("F2").Value = 0.7763925801715814
$ws.Range("L2").Value = 0.2122146180742703
$ws.Range("N2").Value = 1.152311375679417
$ws.Range("O2").Value = 2.661003845049692
$ws.Range("B3").Value = 1.090023758153961
$ws.Range("C3").Value = 0.2070080621704307
$ws.Range("D3").Value = 0.02224389641261837
$ws.Range("E3").Value = 0.09950659437024001
$ws.Range("F3").Value = 0.7690441028765207
$ws.Range("L3").Value = 0.2029099639223659
$ws.Range("N3").Value = 1.163704474325364
$ws.Range("O3").Value = 2.651355265332029
$ws.Range("B4").Value = 1.025899365143459
$ws.Range("C4").Value = 0.2007727506099855
$ws.Range("D4").Value = 0.0216203145965288
$ws.Range("E4").Value = 0.09992322555516786
$ws.Range("F4").Value = 0.7650888328931487
$ws.Range("L4").Value = 0.1973047040672782
$ws.Range("N4").Value = 1.171190191582248
$ws.Range("O4").Value = 2.647313203222978
$ws.Range("B5").Value = 0.9998263544598558
$ws.Range("C5").Value = 0.1982203687486788
$ws.Range("D5").Value = 0.02136507372882335
$ws.Range("E5").Value = 0.1001059769190196
$ws.Range("F5").Value = 0.7636168359998194
$ws.Range("L5").Value = 0.1950476885300958
$ws.Range("N5").Value = 1.174364012167331
$ws.Range("O5").Value = 2.646138631301142
$ws.Range("B6").Value = 0.9955005050966292
$ws.Range("C6").Value = 0.197795861250512
$ws.Range("D6").Value = 0.02132262365556059
$ws.Range("E6").Value = 0.1001371060546834
$ws.Range("F6").Value = 0.7633808506751265
$ws.Range("L6").Value = 0.1946745557268628
$ws.Range("N6").Value = 1.174898472694665
$ws.Range("O6").Value = 2.645972118201001
$ws.Range("B7").Value = 1.02554749799765
$ws.Range("C7").Value = 0.2007383744134899
$ws.Range("D7").Value = 0.02161687686643532
$ws.Range("E7").Value = 0.09992563768544294
$ws.Range("F7").Value = 0.765068415199238
$ws.Range("L7").Value = 0.1972741550658412
$ws.Range("N7").Value = 1.171232495426977
$ws.Range("O7").Value = 2.647295449918374
$ws.Range("B8").Value = 1.158645690242452
$ws.Range("C8").Value = 0.2136221702905345
$ws.Range("D8").Value = 0.02290541910633337
$ws.Range("E8").Value = 0.09910634624975501
$ws.Range("F8").Value = 0.7737431608609384
$ws.Range("L8").Value = 0.2089840230086111
$ws.Range("N8").Value = 1.156137974601954
$ws.Range("O8").Value = 2.65728593382272
$ws.Range("B9").Value = 1.421400947532561
$ws.Range("C9").Value = 0.2385155255674363
$ws.Range("D9").Value = 0.02539554518414633
$ws.Range("E9").Value = 0.09790903914421634
$ws.Range("F9").Value = 0.7951821495015281
$ws.Range("L9").Value = 0.2328015675964537
$ws.Range("N9").Value = 1.130426422940857
$ws.Range("O9").Value = 2.691848309765248
$ws.Range("B10").Value = 1.615474074366148
$ws.Range("C10").Value = 0.2565719894211611
$ws.Range("D10").Value = 0.02720198638466798
$ws.Range("E10").Value = 0.0972790544298423
$ws.Range("F10").Value = 0.8136502656608258
$ws.Range("L10").Value = 0.250821888815608
$ws.Range("N10").Value = 1.113904723711684
$ws.Range("O10").Value = 2.726426362757053
$ws.Range("B11").Value = 1.703978700736741
$ws.Range("C11").Value = 0.2647346361299867
$ws.Range("D11").Value = 0.0280186437467691
$ws.Range("E11").Value = 0.09704671115872188
$ws.Range("F11").Value = 0.8226457974670609
$ws.Range("L11").Value = 0.2591333208486049
$ws.Range("N11").Value = 1.106902547696478
$ws.Range("O11").Value = 2.744164293284797
$ws.Range("B12").Value = 1.7375236501166
$ws.Range("C12").Value = 0.2678181003557825
$ws.Range("D12").Value = 0.02832714209627341
$ws.Range("E12").Value = 0.09696653071955375
$ws.Range("F12").Value = 0.826137886513294
$ws.Range("L12").Value = 0.2622969971709921
$ws.Range("N12").Value = 1.104324859715156
$ws.Range("O12").Value = 2.751170867853517
$ws.Range("B13").Value = 1.730297826171238
$ws.Range("C13").Value = 0.2671543598127073
$ws.Range("D13").Value = 0.02826073519808858
$ws.Range("E13").Value = 0.09698345192268221
$ws.Range("F13").Value = 0.8253819887171829
$ws.Range("L13").Value = 0.2616149181192924
$ws.Range("N13").Value = 1.104876724474586
$ws.Range("O13").Value = 2.749648983206214
$ws.Range("B14").Value = 1.706737867512857
$ws.Range("C14").Value = 0.2649884669417588
$ws.Range("D14").Value = 0.02804403926982246
$ws.Range("E14").Value = 0.09703995826483514
$ws.Range("F14").Value = 0.8229313754082028
$ws.Range("L14").Value = 0.2593932717274043
$ws.Range("N14").Value = 1.106688998846295
$ws.Range("O14").Value = 2.744734919531055
$ws.Range("B15").Value = 1.692310598760685
$ws.Range("C15").Value = 0.2636608063147889
$ws.Range("D15").Value = 0.0279112082690034
$ws.Range("E15").Value = 0.09707558629865787
$ws.Range("F15").Value = 0.8214414680376763
$ws.Range("L15").Value = 0.2580345723623054
$ws.Range("N15").Value = 1.107808691947348
$ws.Range("O15").Value = 2.741762653354158
$ws.Range("B16").Value = 1.609694388645664
$ws.Range("C16").Value = 0.2560374936823848
$ws.Range("D16").Value = 0.02714851166640386
$ws.Range("E16").Value = 0.09729533033493176
$ws.Range("F16").Value = 0.813074363832925
$ws.Range("L16").Value = 0.2502810038860588
$ws.Range("N16").Value = 1.114372669913891
$ws.Range("O16").Value = 2.725307628259515
$ws.Range("B17").Value = 1.55906728881871
$ws.Range("C17").Value = 0.2513475691730207
$ws.Range("D17").Value = 0.02667930216690451
$ws.Range("E17").Value = 0.09744403065681517
$ws.Range("F17").Value = 0.8080937843583627
$ws.Range("L17").Value = 0.2455535692939037
$ws.Range("N17").Value = 1.118531026435583
$ws.Range("O17").Value = 2.715727941592206
$ws.Range("B18").Value = 1.529968693121248
$ws.Range("C18").Value = 0.2486452275189777
$ws.Range("D18").Value = 0.02640894616511247
$ws.Range("E18").Value = 0.09753466427134683
$ws.Range("F18").Value = 0.8052850114918186
$ws.Range("L18").Value = 0.2428451981714375
$ws.Range("N18").Value = 1.120971146377777
$ws.Range("O18").Value = 2.710406931092052
$ws.Range("B19").Value = 1.520120020252534
$ws.Range("C19").Value = 0.2477294382820219
$ws.Range("D19").Value = 0.02631732654351993
$ws.Range("E19").Value = 0.09756622796602699
$ws.Range("F19").Value = 0.8043436077377635
$ws.Range("L19").Value = 0.2419300339967805
$ws.Range("N19").Value = 1.121805631386245
$ws.Range("O19").Value = 2.708637759217993
$ws.Range("B20").Value = 1.564454490506932
$ws.Range("C20").Value = 0.2518473198576032
$ws.Range("D20").Value = 0.02672930003986096
$ws.Range("E20").Value = 0.09742767286929421
$ws.Range("F20").Value = 0.8086181859662247
$ws.Range("L20").Value = 0.2460557033023747
$ws.Range("N20").Value = 1.118083358770988
$ws.Range("O20").Value = 2.716728151698049
$ws.Range("B21").Value = 1.713657191149025
$ws.Range("C21").Value = 0.2656248485309902
$ws.Range("D21").Value = 0.02810770867422718
$ws.Range("E21").Value = 0.09702314918738786
$ws.Range("F21").Value = 0.8236488529773851
$ws.Range("L21").Value = 0.2600453807628043
$ws.Range("N21").Value = 1.10615468451509
$ws.Range("O21").Value = 2.746170432751711
$ws.Range("B22").Value = 1.81134496030586
$ws.Range("C22").Value = 0.2745851564137638
$ws.Range("D22").Value = 0.02900418625712575
$ws.Range("E22").Value = 0.09680425247724855
$ws.Range("F22").Value = 0.8339717507095372
$ws.Range("L22").Value = 0.2692835256877402
$ws.Range("N22").Value = 1.098789274595958
$ws.Range("O22").Value = 2.76710108996545
$ws.Range("B23").Value = 1.759191604286855
$ws.Range("C23").Value = 0.2698069635565616
$ws.Range("D23").Value = 0.02852612725332904
$ws.Range("E23").Value = 0.09691691899238997
$ws.Range("F23").Value = 0.8284164510124015
$ws.Range("L23").Value = 0.2643442739981481
$ws.Range("N23").Value = 1.102680916745456
$ws.Range("O23").Value = 2.755775239953635
$ws.Range("B24").Value = 1.562018912846838
$ws.Range("C24").Value = 0.2516214013804756
$ws.Range("D24").Value = 0.02670669787294599
$ws.Range("E24").Value = 0.0974350522029539
$ws.Range("F24").Value = 0.8083809338498043
$ws.Range("L24").Value = 0.2458286589417753
$ws.Range("N24").Value = 1.11828559536287
$ws.Range("O24").Value = 2.716275375893275
$ws.Range("B25").Value = 1.350135486689737
$ws.Range("C25").Value = 0.2318216392692705
$ws.Range("D25").Value = 0.02472589985626428
$ws.Range("E25").Value = 0.0981891070961467
$ws.Range("F25").Value = 0.7889064867293456
$ws.Range("L25").Value = 0.2262667429534702
$ws.Range("N25").Value = 1.136965999057317
$ws.Range("O25").Value = 2.680889486086016

Write-Output "Updated 192 cells"
